# Update of the GL (Group Learning) results
# - sheet "5vs8": insert two columns (J:K) to make room for a new "default"
#   hyper-parameter result column; fill the new J (Training "default") and
#   T (Test "default") columns; add a STD row (row 25) under the Average row.
# - sheet "F27vsF33": add a STD row (row 25) under the Average row.
# - minor selection/view tweaks.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Sheet "5vs8" (sheet1.xml)
# ---------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("5vs8")

# Insert two blank columns before the Test_error c-sweep block (J:K).
# This shifts the existing J:Q (Test_error c=0.001..10000) block to L:S,
# carrying formulas, values and styles (incl. the yellow-highlighted
# "best c" column) along with it.
$ws1.Columns("J:K").Insert()

# New "default" header cells (shared string "default")
$ws1.Range("J3").Value = "default"
$ws1.Range("T3").Value = "default"

# New "default" Training_error column (all zeros)
for ($r = 4; $r -le 23; $r++) {
    $ws1.Cells.Item($r, 10).Value = 0
}

# New "default" Test_error column values
$ws1.Cells.Item(4, 20).Value = 0.211
$ws1.Cells.Item(5, 20).Value = 0.149
$ws1.Cells.Item(6, 20).Value = 0.116
$ws1.Cells.Item(7, 20).Value = 0.196
$ws1.Cells.Item(8, 20).Value = 0.212
$ws1.Cells.Item(9, 20).Value = 0.148
$ws1.Cells.Item(10, 20).Value = 0.301
$ws1.Cells.Item(11, 20).Value = 0.279
$ws1.Cells.Item(12, 20).Value = 0.182
$ws1.Cells.Item(13, 20).Value = 0.148
$ws1.Cells.Item(14, 20).Value = 0.181
$ws1.Cells.Item(15, 20).Value = 0.299
$ws1.Cells.Item(16, 20).Value = 0.205
$ws1.Cells.Item(17, 20).Value = 0.268
$ws1.Cells.Item(18, 20).Value = 0.204
$ws1.Cells.Item(19, 20).Value = 0.172
$ws1.Cells.Item(20, 20).Value = 0.241
$ws1.Cells.Item(21, 20).Value = 0.192
$ws1.Cells.Item(22, 20).Value = 0.276
$ws1.Cells.Item(23, 20).Value = 0.184

# Averages for the two new "default" columns (row 24)
$ws1.Range("J24").Formula = "=AVERAGE(J4:J23)"
$ws1.Range("T24").Formula = "=AVERAGE(T4:T23)"

# New STD row (row 25): label + STDEV for every data column (K and the
# spacer column to the right of the "default" Training column stay blank)
$ws1.Range("A25").Value = "STD"
$stdCols1 = @("B","C","D","E","F","G","H","I","J","L","M","N","O","P","Q","R","S","T")
foreach ($col in $stdCols1) {
    $ws1.Range("$col 25".Replace(" ", "")).Formula = "=STDEV(" + $col + "4:" + $col + "23)"
}

# View tweaks on sheet "5vs8"
$ws1.Range("V14").Select() | Out-Null

# ---------------------------------------------------------------------
# Sheet "F27vsF33" (sheet7.xml)
# ---------------------------------------------------------------------
$ws7 = $wb.Worksheets.Item("F27vsF33")

$stdCols7 = @("B","C","D","E","F","G","H","I","J","K","L","M","N","O","P","Q")
foreach ($col in $stdCols7) {
    $ws7.Range("$col 25".Replace(" ", "")).Formula = "=STDEV(" + $col + "4:" + $col + "23)"
}

# Re-select sheet "5vs8" as the active sheet/tab (it was tabSelected="1" before)
$ws1.Activate() | Out-Null
